$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '84.280.13'
$ws.Range('E2').Value = '  +5.65%  '
$ws.Range('D3').Value = '3.284.01'
$ws.Range('E3').Value = '  +2.13%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.78'
$ws.Range('E5').Value = '  +4.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '638.51'
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.321'
$ws.Range('E7').Value = '  +21.42%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.590'
$ws.Range('E9').Value = '  -1.82%  '
$ws.Range('D10').Value = '3.289.47'
$ws.Range('E10').Value = '  +2.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.599'
$ws.Range('E11').Value = '  -0.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000278'
$ws.Range('E12').Value = '  +2.51%  '
$ws.Range('E13').Value = '  -0.05%  '
$ws.Range('D14').Value = '3.887.90'
$ws.Range('E14').Value = '  +2.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.42'
$ws.Range('E15').Value = '  -0.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '33.62'
$ws.Range('E16').Value = '  +2.79%  '
$ws.Range('D17').Value = '84.352.47'
$ws.Range('E17').Value = '  +5.91%  '
$ws.Range('D18').Value = '3.281.79'
$ws.Range('E18').Value = '  +2.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.20'
$ws.Range('E19').Value = '  +5.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.47'
$ws.Range('E20').Value = '  -1.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '450.35'
$ws.Range('E21').Value = '  +0.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.15'
$ws.Range('E22').Value = '  -3.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.25'
$ws.Range('E23').Value = '  -0.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.46'
$ws.Range('E24').Value = '  +5.88%  '
$ws.Range('B25').Value = 'Aptos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.47'
$ws.Range('E25').Value = '  +13.92%  '
$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.27'
$ws.Range('E26').Value = '  +8.75%  '
$ws.Range('D27').Value = '3.447.33'
$ws.Range('E27').Value = '  +2.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '78.15'
$ws.Range('E28').Value = '  +0.33%  '
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0000126'
$ws.Range('E30').Value = '  -0.90%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.22'
$ws.Range('E31').Value = '  -0.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.993'
$ws.Range('E32').Value = '  -0.46%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '571.84'
$ws.Range('E33').Value = '  +1.23%  '
$ws.Range('B34').Value = 'Cronos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.154'
$ws.Range('E34').Value = '  +26.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.51'
$ws.Range('E35').Value = '  -1.28%  '
$ws.Range('E36').Value = '  -2.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.03'
$ws.Range('E37').Value = '  -1.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '23.29'
$ws.Range('E38').Value = '  +0.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.24'
$ws.Range('E39').Value = '  +8.56%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.414'
$ws.Range('E41').Value = '  -0.23%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.05'
$ws.Range('E42').Value = '  +12.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.92'
$ws.Range('E43').Value = '  +3.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.07'
$ws.Range('E44').Value = '  +12.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '159.79'
$ws.Range('E45').Value = '  -2.35%  '
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '191.63'
$ws.Range('E47').Value = '  -1.75%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '45.16'
$ws.Range('E48').Value = '  +4.95%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.34'
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.781'
$ws.Range('E50').Value = '  -2.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '26.53'
$ws.Range('E51').Value = '  +2.16%  '
